$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the date-range header text in D3
#    "01/01/2024 to 31/02/2025" -> "01/01/2024 to 31/01/2025"
# ------------------------------------------------------------------
$ws.Range("D3").Value() = "01/01/2024 to 31/01/2025"

# ------------------------------------------------------------------
# 2. Prepare rows 12 and 13 by copying the formatting of row 11
#    (only the B, D, G, H columns are populated, matching the
#    existing transaction rows 9-11)
# ------------------------------------------------------------------
$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B13").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D13").PasteSpecial(-4122)

$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$ws.Range("G13").PasteSpecial(-4122)

$ws.Range("H11").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H13").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3. Row 12 becomes what row 11 used to be: the Phone Company
#    direct debit on 31/01/2024 (serial 45322)
# ------------------------------------------------------------------
$ws.Range("B12").Value() = 45322
$ws.Range("D12").Value() = "DIRECT DEBIT PAYMENT TO Phone Company"
$ws.Range("G12").Value() = 10
$ws.Range("H12").Formula = "=SUM(H10,F12,-G12)"

# ------------------------------------------------------------------
# 4. Row 11 is replaced with a new transaction: Netflix card
#    payment on 05/01/2024 (serial 45296)
# ------------------------------------------------------------------
$ws.Range("B11").Value() = 45296
$ws.Range("D11").Value() = "CARD PAYMENT TO Netflix ON 05-01-2024"

# ------------------------------------------------------------------
# 5. Row 13 is a brand new transaction: Rent direct debit
# ------------------------------------------------------------------
$ws.Range("B13").Value() = 45322
$ws.Range("D13").Value() = "DIRECT DEBIT PAYMENT TO Rent"
$ws.Range("G13").Value() = 450
$ws.Range("H13").Formula = "=SUM(H11,F13,-G13)"

# ------------------------------------------------------------------
# 6. Update the active selection to H19
# ------------------------------------------------------------------
$null = $ws.Range("H19").Select()
